$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.441.06"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "1.805.92"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'1.005"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'307.09"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("D7").Value = "'0.4521"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("D8").Value = "'0.3598"
$ws.Range("E8").Value = "  -1.98%  "
$ws.Range("D9").Value = "'46.41"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'0.07088"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.8900"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "'0.07804"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'19.47"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "1.822.79"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'5.290"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'6.336"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "'85.10"
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D19").Value = "'0.000008475"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "26.474.82"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "'14.27"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'4.971"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "2.040.98"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value = "'10.52"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'1.962"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").Value = "'150.94"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'17.82"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "'2.053"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").Value = "'111.98"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "'4.864"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'0.08702"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").Value = "'3.117"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'2.845"
$ws.Range("E34").Value = "  +14.61%  "
$ws.Range("D35").Value = "'4.444"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'0.7212"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").Value = "'1.111"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "'0.01940"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "'0.05099"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'2.895"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'0.5123"
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("D44").Value = "'6.783"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "'0.1513"
$ws.Range("E45").Value = "  -4.79%  "
$ws.Range("D46").Value = "'8.018"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").Value = "'0.4658"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "'10.02"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'100.37"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("D51").Value = "'1.573"
$ws.Range("E51").Value = "  -1.01%  "
